$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range('D2').Value = '29.477.08'
$ws.Range('D3').Value = '1.879.23'
$ws.Range('E3').Value = '  +1.27%  '
Set-TextValue 'D4' '0.9990'
$ws.Range('E4').Value = '  -0.12%  '
Set-TextValue 'D5' '0.7160'
$ws.Range('E5').Value = '  +1.74%  '
Set-TextValue 'D6' '242.21'
$ws.Range('E6').Value = '  +1.86%  '
Set-TextValue 'D7' '0.9999'
$ws.Range('E7').Value = '  -0.06%  '
Set-TextValue 'D8' '0.07893'
$ws.Range('E8').Value = '  -1.37%  '
Set-TextValue 'D9' '0.3126'
$ws.Range('E9').Value = '  +3.40%  '
Set-TextValue 'D10' '25.36'
$ws.Range('E10').Value = '  +7.47%  '
Set-TextValue 'D11' '0.08276'
$ws.Range('E11').Value = '  +0.94%  '
$ws.Range('D12').Value = '1.930.81'
$ws.Range('E12').Value = '  +4.94%  '
Set-TextValue 'D13' '0.7323'
$ws.Range('E13').Value = '  +3.71%  '
Set-TextValue 'D14' '5.298'
$ws.Range('E14').Value = '  +2.10%  '
Set-TextValue 'D15' '91.42'
$ws.Range('E15').Value = '  +2.05%  '
$ws.Range('D16').Value = '29.510.42'
$ws.Range('E16').Value = '  +1.14%  '
Set-TextValue 'D17' '5.957'
$ws.Range('E17').Value = '  +2.31%  '
Set-TextValue 'D18' '247.80'
$ws.Range('E18').Value = '  +4.74%  '
Set-TextValue 'D19' '0.000007880'
$ws.Range('E19').Value = '  +0.42%  '
Set-TextValue 'D20' '13.37'
$ws.Range('E20').Value = '  +1.30%  '
Set-TextValue 'D21' '0.9991'
$ws.Range('E21').Value = '  -0.08%  '
Set-TextValue 'D22' '8.012'
$ws.Range('E22').Value = '  +6.74%  '
Set-TextValue 'D23' '0.9991'
$ws.Range('E23').Value = '  -0.12%  '
Set-TextValue 'D24' '0.1610'
$ws.Range('E24').Value = '  +13.65%  '
Set-TextValue 'D25' '163.94'
$ws.Range('E25').Value = '  +0.61%  '
Set-TextValue 'D26' '9.059'
$ws.Range('E26').Value = '  +2.08%  '
Set-TextValue 'D27' '18.38'
$ws.Range('E27').Value = '  +1.59%  '
Set-TextValue 'D28' '1.362'
$ws.Range('E28').Value = '  -2.72%  '
Set-TextValue 'D29' '1.498'
$ws.Range('E29').Value = '  +1.66%  '
Set-TextValue 'D30' '4.394'
$ws.Range('E30').Value = '  +1.26%  '
Set-TextValue 'D31' '4.126'
$ws.Range('E31').Value = '  +2.67%  '
Set-TextValue 'D32' '0.05294'
$ws.Range('E32').Value = '  +2.43%  '
Set-TextValue 'D33' '1.958'
$ws.Range('E33').Value = '  +2.48%  '
$ws.Range('E34').Value = '  +2.96%  '
Set-TextValue 'D35' '0.7288'
$ws.Range('E35').Value = '  +2.45%  '
Set-TextValue 'D36' '2.674'
$ws.Range('E36').Value = '  -0.20%  '
Set-TextValue 'D37' '0.01874'
$ws.Range('E37').Value = '  +1.41%  '
$ws.Range('D38').Value = '1.225.92'
$ws.Range('E38').Value = '  +6.08%  '
$ws.Range('E39').Value = '  +0.88%  '
$ws.Range('B40').Value = 'Aave'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D40' '75.09'
$ws.Range('E40').Value = '  +6.99%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D41' '0.9126'
$ws.Range('E41').Value = '  -1.87%  '
Set-TextValue 'D42' '6.211'
$ws.Range('E42').Value = '  +3.83%  '
Set-TextValue 'D43' '1.000'
$ws.Range('E43').Value = '  -0.04%  '
Set-TextValue 'D44' '102.89'
$ws.Range('E44').Value = '  +0.19%  '
$ws.Range('D45').Value = '2.043.18'
$ws.Range('E45').Value = '  +2.85%  '
Set-TextValue 'D46' '0.5264'
$ws.Range('E46').Value = '  -0.23%  '
Set-TextValue 'D47' '2.991'
$ws.Range('E47').Value = '  +15.53%  '
Set-TextValue 'D48' '1.780'
$ws.Range('E48').Value = '  +2.19%  '
Set-TextValue 'D49' '9.357'
$ws.Range('E49').Value = '  +2.08%  '
Set-TextValue 'D50' '0.4343'
$ws.Range('E50').Value = '  +2.02%  '
Set-TextValue 'D51' '7.130'
$ws.Range('E51').Value = '  +2.21%  '
